# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook's "K" column (column G) held values that were computed from the
# old "Strike#" metric. This re-generates those values (the new K counts) for
# every data row (rows 2-65) and writes them back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..65 (one value per row, in row order).
$kValues = @(
    0, 2, 2, 2, 2, 1, 1, 2, 0, 1,
    0, 2, 1, 3, 3, 0, 1, 1, 1, 0,
    0, 1, 1, 3, 1, 1, 2, 1, 0, 1,
    1, 1, 1, 1, 2, 2, 2, 3, 1, 0,
    0, 0, 2, 4, 0, 2, 3, 2, 1, 1,
    0, 1, 2, 4, 3, 1, 6, 2, 5, 4,
    3, 1, 2, 1
)

$startRow = 2
$endRow = 65

# Write the value for each row individually into column G (7).
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
